$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Item ID 3) was missing its item name -> set to "Burger"
$ws.Range("B4").Value = "Burger"

# Row 5 (Item ID 4, Pizza) "Special Item" flag was off -> mark it a special item
$ws.Range("D5").Value = $true

# Move the active selection to B5 (matches the saved selection state)
$ws.Range("B5").Select()
